# Auto-generated edit script: updates cryptos price/volume table
# to match the Fri Mar 15 06:34:38 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.405.07"
$ws.Range("E2").Value = "  -6.94%  "

# Row 3
$ws.Range("D3").Value = "3.737.04"
$ws.Range("E3").Value = "  -6.26%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.95%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.49%  "

# Row 7
$ws.Range("D7").Value = "3.876.80"
$ws.Range("E7").Value = "  -2.42%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.637"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.55%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.723"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.89%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.43%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.88%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000304"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.86%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.63%  "

# Row 15
$ws.Range("D15").Value = "4.297.48"
$ws.Range("E15").Value = "  -7.10%  "

# Row 16
$ws.Range("D16").Value = "3.736.05"
$ws.Range("E16").Value = "  -6.26%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.56"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.58%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.72%  "

# Row 19
$ws.Range("E19").Value = "  -7.49%  "

# Row 20
$ws.Range("E20").Value = "  -2.98%  "

# Row 21
$ws.Range("D21").Value = "68.176.68"
$ws.Range("E21").Value = "  -7.01%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "412.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.99%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.81%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.54%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.49%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.58%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.72%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.08%  "

# Row 33
$ws.Range("E33").Value = "  -6.80%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.119"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.46%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.30%  "

# Row 36
$ws.Range("D36").Value = "0.0₃0943"
$ws.Range("E36").Value = "  -9.59%  "

# Row 37
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "617.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.34%  "

# Row 38
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "65.57"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.56%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.407"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.67%  "

# Row 40
$ws.Range("E40").Value = "  +0.24%  "

# Row 41
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.40%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.995"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.139"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.40%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.92%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0445"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.33%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.80%  "

# Row 48
$ws.Range("E48").Value = "  -8.71%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -14.35%  "

# Row 50
$ws.Range("D50").Value = "2.748.64"
$ws.Range("E50").Value = "  -2.89%  "

# Row 51
$ws.Range("E51").Value = "  -7.01%  "

